$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VLO")

$values = @{
    "B2"  = 2254000000.0
    "B3"  = 9014000000.0
    "B7"  = 30379000000.0
    "B11" = 36023000000.0
    "B12" = 53614000000.0
    "B13" = 9113000000.0
    "B14" = 1087000000.0
    "B16" = 1373000000.0
    "B21" = 5034000000.0
    "B22" = 3616000000.0
    "B23" = 22580000000.0
    "B25" = 6810000000.0
    "B26" = 7000000.0
    "B27" = 27849000000.0
    "B28" = 15700000000.0
    "B31" = 53614000000.0
    "B32" = 408761000.0
    "B33" = 18727000000.0
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
